# preparation publication 0.2.0
# - bump Version to 0.2.0
# - bump Date to the new publication timestamp
# - insert a new "Jurisdiction" / "iso:code:3166:FR" row right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "0.2.0"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Insert a new row right after "Contact" (row 10) so "Jurisdiction" lands on row 11,
# pushing "Description" and everything below it down by one row.
$newRow = $ws.Rows.Item(11)
$newRow.Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Match the formatting used by the other data rows (the row below, which still carries
# the original data-row style) instead of the bare default the insert produced.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
